# Updates as of 28th April 2020
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 40: fill in the previously-blank "Travelled From" cell ---
$ws.Cells.Item(40, 4).Value = "None"

# --- New data rows 43-46 ---
# New shared-string-valued cells are written first, in the exact order the
# source workbook introduced them, so the shared-string table grows with the
# same new unique entries (110-116) as the authored edit.
$ws.Cells.Item(46, 5).Value = "Mombasa(4),Nairobi(7)"
$ws.Cells.Item(46, 12).Value = "3-75."
$ws.Cells.Item(44, 12).Value = "14-60"
$ws.Cells.Item(44, 5).Value = "Mombasa(4),Nairobi(8)"
$ws.Cells.Item(44, 7).Value = "Community(12)"
$ws.Cells.Item(45, 12).Value = "13-65"
$ws.Cells.Item(45, 5).Value = "Mombasa(4),Nairobi(3),Kwale(1)"

# Row 43 (25-Apr-2020)
$ws.Cells.Item(43, 1).Value = 43946
$ws.Cells.Item(43, 1).NumberFormat = "d-mmm-yy"
$ws.Cells.Item(43, 1).HorizontalAlignment = -4108
$ws.Cells.Item(43, 2).Value = 7
$ws.Cells.Item(43, 4).Value = "None"
$ws.Cells.Item(43, 6).Value = 343
$ws.Cells.Item(43, 7).Value = "Community(7)"
$ws.Cells.Item(43, 8).Value = 0
$ws.Cells.Item(43, 9).Value = 0

# Row 44 (26-Apr-2020)
$ws.Cells.Item(44, 1).Value = 43947
$ws.Cells.Item(44, 1).NumberFormat = "d-mmm-yy"
$ws.Cells.Item(44, 1).HorizontalAlignment = -4108
$ws.Cells.Item(44, 2).Value = 12
$ws.Cells.Item(44, 4).Value = "None"
$ws.Cells.Item(44, 6).Value = 355
$ws.Cells.Item(44, 8).Value = 8
$ws.Cells.Item(44, 9).Value = 0

# Row 45 (27-Apr-2020)
$ws.Cells.Item(45, 1).Value = 43948
$ws.Cells.Item(45, 1).NumberFormat = "d-mmm-yy"
$ws.Cells.Item(45, 1).HorizontalAlignment = -4108
$ws.Cells.Item(45, 2).Value = 8
$ws.Cells.Item(45, 3).Value = 366
$ws.Cells.Item(45, 4).Value = "None"
$ws.Cells.Item(45, 6).Value = 363
$ws.Cells.Item(45, 7).Value = "Community(8)"
$ws.Cells.Item(45, 8).Value = 0
$ws.Cells.Item(45, 9).Value = 0

# Row 46 (28-Apr-2020)
$ws.Cells.Item(46, 1).Value = 43949
$ws.Cells.Item(46, 1).NumberFormat = "d-mmm-yy"
$ws.Cells.Item(46, 1).HorizontalAlignment = -4108
$ws.Cells.Item(46, 2).Value = 11
$ws.Cells.Item(46, 3).Value = 579
$ws.Cells.Item(46, 4).Value = "None"
$ws.Cells.Item(46, 6).Value = 374
$ws.Cells.Item(46, 7).Value = "Community(11)"
$ws.Cells.Item(46, 8).Value = 10
$ws.Cells.Item(46, 9).Value = 0

# --- View state: selection left on G44 (matches the authored edit) ---
$null = $ws.Range("G44").Select()
